{"js": "// Replace the date and the 25 two-digit multiplication problems with\n// their new values, as described by the diff. Every old value is\n// unique in the document, so a simple exact, case-sensitive search +\n// replace per value is safe and unambiguous.\nconst replacements = [\n  [\"2025-06-03 Tuesday\", \"2025-06-04 Wednesday\"],\n  [\"76\u00d746=\", \"37\u00d773=\"],\n  [\"53\u00d767=\", \"74\u00d744=\"],\n  [\"39\u00d794=\", \"60\u00d728=\"],\n  [\"50\u00d764=\", \"26\u00d713=\"],\n  [\"51\u00d751=\", \"32\u00d777=\"],\n  [\"76\u00d767=\", \"50\u00d772=\"],\n  [\"50\u00d767=\", \"41\u00d786=\"],\n  [\"89\u00d769=\", \"53\u00d725=\"],\n  [\"70\u00d770=\", \"98\u00d779=\"],\n  [\"64\u00d712=\", \"19\u00d780=\"],\n  [\"26\u00d790=\", \"96\u00d781=\"],\n  [\"26\u00d745=\", \"20\u00d737=\"],\n  [\"54\u00d768=\", \"76\u00d742=\"],\n  [\"68\u00d731=\", \"67\u00d725=\"],\n  [\"11\u00d795=\", \"88\u00d717=\"],\n  [\"36\u00d779=\", \"72\u00d754=\"],\n  [\"63\u00d772=\", \"95\u00d763=\"],\n  [\"28\u00d717=\", \"14\u00d759=\"],\n  [\"81\u00d740=\", \"28\u00d765=\"],\n  [\"12\u00d793=\", \"81\u00d792=\"],\n  [\"37\u00d714=\", \"36\u00d750=\"],\n  [\"76\u00d769=\", \"26\u00d767=\"],\n  [\"70\u00d798=\", \"93\u00d798=\"],\n  [\"34\u00d740=\", \"13\u00d783=\"],\n  [\"99\u00d730=\", \"69\u00d713=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and the 25 two-digit multiplication problems with\n# their new values, as described by the diff. Every old value is\n# unique in the document, so a simple exact, case-sensitive\n# Find/Replace per value is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-06-03 Tuesday\", \"2025-06-04 Wednesday\"),\n  @(\"76\u00d746=\", \"37\u00d773=\"),\n  @(\"53\u00d767=\", \"74\u00d744=\"),\n  @(\"39\u00d794=\", \"60\u00d728=\"),\n  @(\"50\u00d764=\", \"26\u00d713=\"),\n  @(\"51\u00d751=\", \"32\u00d777=\"),\n  @(\"76\u00d767=\", \"50\u00d772=\"),\n  @(\"50\u00d767=\", \"41\u00d786=\"),\n  @(\"89\u00d769=\", \"53\u00d725=\"),\n  @(\"70\u00d770=\", \"98\u00d779=\"),\n  @(\"64\u00d712=\", \"19\u00d780=\"),\n  @(\"26\u00d790=\", \"96\u00d781=\"),\n  @(\"26\u00d745=\", \"20\u00d737=\"),\n  @(\"54\u00d768=\", \"76\u00d742=\"),\n  @(\"68\u00d731=\", \"67\u00d725=\"),\n  @(\"11\u00d795=\", \"88\u00d717=\"),\n  @(\"36\u00d779=\", \"72\u00d754=\"),\n  @(\"63\u00d772=\", \"95\u00d763=\"),\n  @(\"28\u00d717=\", \"14\u00d759=\"),\n  @(\"81\u00d740=\", \"28\u00d765=\"),\n  @(\"12\u00d793=\", \"81\u00d792=\"),\n  @(\"37\u00d714=\", \"36\u00d750=\"),\n  @(\"76\u00d769=\", \"26\u00d767=\"),\n  @(\"70\u00d798=\", \"93\u00d798=\"),\n  @(\"34\u00d740=\", \"13\u00d783=\"),\n  @(\"99\u00d730=\", \"69\u00d713=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #          MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n  #          Format, ReplaceWith, Replace)\n  # Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll)\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
